$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Lrpap1"
$ws.Cells.Item(2,3).Value = "Vldlr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 6.540008666666668
$ws.Cells.Item(2,8).Value = 19.620026
$ws.Cells.Item(2,9).Value = 0.2365207520404831
$ws.Cells.Item(2,10).Value = 0.2365207520404831
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.3374003333333333
$ws.Cells.Item(2,14).Value = 1.012201
$ws.Cells.Item(2,15).Value = 0.01738364872808817
$ws.Cells.Item(2,16).Value = 0.01738364872808818
$ws.Cells.Item(2,17).Value = 2.206601104136222
$ws.Cells.Item(2,18).Value = 19.859409937226
$ws.Cells.Item(2,19).Value = 0.004111593670375002
$ws.Cells.Item(2,20).Value = 0.004111593670375003

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Lrpap1"
$ws.Cells.Item(3,3).Value = "Vldlr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 6.540008666666668
$ws.Cells.Item(3,8).Value = 19.620026
$ws.Cells.Item(3,9).Value = 0.2365207520404831
$ws.Cells.Item(3,10).Value = 0.2365207520404831
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 16.71131166666667
$ws.Cells.Item(3,14).Value = 50.133935
$ws.Cells.Item(3,15).Value = 0.8610055862391021
$ws.Cells.Item(3,16).Value = 0.8610055862391023
$ws.Cells.Item(3,17).Value = 109.2921231313678
$ws.Cells.Item(3,18).Value = 983.6291081823101
$ws.Cells.Item(3,19).Value = 0.2036456887683294
$ws.Cells.Item(3,20).Value = 0.2036456887683295

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Lrpap1"
$ws.Cells.Item(4,3).Value = "Vldlr"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 6.540008666666668
$ws.Cells.Item(4,8).Value = 19.620026
$ws.Cells.Item(4,9).Value = 0.2365207520404831
$ws.Cells.Item(4,10).Value = 0.2365207520404831
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.360351
$ws.Cells.Item(4,14).Value = 7.081053000000001
$ws.Cells.Item(4,15).Value = 0.1216107650328097
$ws.Cells.Item(4,16).Value = 0.1216107650328097
$ws.Cells.Item(4,17).Value = 15.43671599637534
$ws.Cells.Item(4,18).Value = 138.930443967378
$ws.Cells.Item(4,19).Value = 0.02876346960177862
$ws.Cells.Item(4,20).Value = 0.02876346960177863

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Lrpap1"
$ws.Cells.Item(5,3).Value = "Vldlr"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 14.67485766666667
$ws.Cells.Item(5,8).Value = 44.024573
$ws.Cells.Item(5,9).Value = 0.5307192311682535
$ws.Cells.Item(5,10).Value = 0.5307192311682536
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.3374003333333333
$ws.Cells.Item(5,14).Value = 1.012201
$ws.Cells.Item(5,15).Value = 0.01738364872808817
$ws.Cells.Item(5,16).Value = 0.01738364872808818
$ws.Cells.Item(5,17).Value = 4.951301868352555
$ws.Cells.Item(5,18).Value = 44.561716815173
$ws.Cells.Item(5,19).Value = 0.009225836687869943
$ws.Cells.Item(5,20).Value = 0.009225836687869948

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Lrpap1"
$ws.Cells.Item(6,3).Value = "Vldlr"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 14.67485766666667
$ws.Cells.Item(6,8).Value = 44.024573
$ws.Cells.Item(6,9).Value = 0.5307192311682535
$ws.Cells.Item(6,10).Value = 0.5307192311682536
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 16.71131166666667
$ws.Cells.Item(6,14).Value = 50.133935
$ws.Cells.Item(6,15).Value = 0.8610055862391021
$ws.Cells.Item(6,16).Value = 0.8610055862391023
$ws.Cells.Item(6,17).Value = 245.2361201316395
$ws.Cells.Item(6,18).Value = 2207.125081184755
$ws.Cells.Item(6,19).Value = 0.4569522227603876
$ws.Cells.Item(6,20).Value = 0.4569522227603878

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Lrpap1"
$ws.Cells.Item(7,3).Value = "Vldlr"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 14.67485766666667
$ws.Cells.Item(7,8).Value = 44.024573
$ws.Cells.Item(7,9).Value = 0.5307192311682535
$ws.Cells.Item(7,10).Value = 0.5307192311682536
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.360351
$ws.Cells.Item(7,14).Value = 7.081053000000001
$ws.Cells.Item(7,15).Value = 0.1216107650328097
$ws.Cells.Item(7,16).Value = 0.1216107650328097
$ws.Cells.Item(7,17).Value = 34.63781496837434
$ws.Cells.Item(7,18).Value = 311.740334715369
$ws.Cells.Item(7,19).Value = 0.06454117171999586
$ws.Cells.Item(7,20).Value = 0.0645411717199959

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Lrpap1"
$ws.Cells.Item(8,3).Value = "Vldlr"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 6.436020999999999
$ws.Cells.Item(8,8).Value = 19.308063
$ws.Cells.Item(8,9).Value = 0.2327600167912634
$ws.Cells.Item(8,10).Value = 0.2327600167912634
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.3374003333333333
$ws.Cells.Item(8,14).Value = 1.012201
$ws.Cells.Item(8,15).Value = 0.01738364872808817
$ws.Cells.Item(8,16).Value = 0.01738364872808818
$ws.Cells.Item(8,17).Value = 2.171515630740333
$ws.Cells.Item(8,18).Value = 19.543640676663
$ws.Cells.Item(8,19).Value = 0.004046218369843228
$ws.Cells.Item(8,20).Value = 0.004046218369843229

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Lrpap1"
$ws.Cells.Item(9,3).Value = "Vldlr"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 6.436020999999999
$ws.Cells.Item(9,8).Value = 19.308063
$ws.Cells.Item(9,9).Value = 0.2327600167912634
$ws.Cells.Item(9,10).Value = 0.2327600167912634
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 16.71131166666667
$ws.Cells.Item(9,14).Value = 50.133935
$ws.Cells.Item(9,15).Value = 0.8610055862391021
$ws.Cells.Item(9,16).Value = 0.8610055862391023
$ws.Cells.Item(9,17).Value = 107.5543528242117
$ws.Cells.Item(9,18).Value = 967.9891754179049
$ws.Cells.Item(9,19).Value = 0.200407674710385
$ws.Cells.Item(9,20).Value = 0.200407674710385

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Lrpap1"
$ws.Cells.Item(10,3).Value = "Vldlr"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 6.436020999999999
$ws.Cells.Item(10,8).Value = 19.308063
$ws.Cells.Item(10,9).Value = 0.2327600167912634
$ws.Cells.Item(10,10).Value = 0.2327600167912634
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.360351
$ws.Cells.Item(10,14).Value = 7.081053000000001
$ws.Cells.Item(10,15).Value = 0.1216107650328097
$ws.Cells.Item(10,16).Value = 0.1216107650328097
$ws.Cells.Item(10,17).Value = 15.191268603371
$ws.Cells.Item(10,18).Value = 136.721417430339
$ws.Cells.Item(10,19).Value = 0.02830612371103516
$ws.Cells.Item(10,20).Value = 0.02830612371103517
